$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.849.13'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').Value = '2.661.67'
$ws.Range('E3').Value = '  +1.45%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '609.23'
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('D6').Value = '157.55'
$ws.Range('E6').Value = '  +3.40%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '0.125'
$ws.Range('E9').Value = '  +9.10%  '
$ws.Range('D10').Value = '6.04'
$ws.Range('E10').Value = '  +4.09%  '
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '0.0000208'
$ws.Range('E13').Value = '  +21.18%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '30.09'
$ws.Range('E14').Value = '  +5.06%  '
$ws.Range('D15').Value = '3.139.88'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '65.682.96'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '2.655.20'
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').Value = '12.67'
$ws.Range('E18').Value = '  +3.24%  '
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').Value = '359.30'
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('D21').Value = '7.46'
$ws.Range('E21').Value = '  +4.00%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '70.19'
$ws.Range('E23').Value = '  +3.71%  '
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000107'
$ws.Range('E25').Value = '  +17.27%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '9.51'
$ws.Range('E26').Value = '  +2.24%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.172'
$ws.Range('E27').Value = '  +4.80%  '
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').Value = '1.64'
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').Value = '8.11'
$ws.Range('E29').Value = '  -1.67%  '
$ws.Range('D30').Value = '2.20'
$ws.Range('E30').Value = '  +6.60%  '
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').Value = '537.75'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').Value = '5.54'
$ws.Range('E34').Value = '  -3.11%  '
$ws.Range('D35').Value = '6.45'
$ws.Range('E35').Value = '  +3.67%  '
$ws.Range('D36').Value = '0.433'
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('D37').Value = '20.77'
$ws.Range('E37').Value = '  +3.33%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '2.02'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '162.76'
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = '169.81'
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('D43').Value = '42.10'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').Value = '4.17'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').Value = '2.34'
$ws.Range('E45').Value = '  +5.06%  '
$ws.Range('D46').Value = '0.0612'
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('D47').Value = '23.09'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').Value = '0.659'
$ws.Range('E48').Value = '  +2.82%  '
$ws.Range('D49').Value = '0.0265'
$ws.Range('E49').Value = '  +5.55%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '19.87'
$ws.Range('E50').Value = '  +2.82%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.0979'
$ws.Range('E51').Value = '  -0.34%  '
